$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (A3/B3) that was appended to the sheet.
$ws.Range("A3").Value = "This is just an example"
$ws.Range("B3").Value = 1

# Move the selection down to A4, as it was left after entering the new row.
$ws.Range("A4").Select() | Out-Null

# The workbook's single built-in cell style was relabelled from the
# Dutch default name "Standaard" to the English default name "Normal".
# The engine doesn't support renaming a style in place, so recreate it
# under the new name (keeps count at 1 / builtinId 0, same as Excel would).
$wb.Styles.Item("Standaard").Delete()
$wb.Styles.Add("Normal")
